$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/percentage strings stay as text (matches source data which
# stores these as plain text, not numbers), without leaving a residual style on any cell.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "22.336.39"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "1.562.69"
$ws.Range("E3").Value = "  -0.92%  "

$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").Value = "0.9996"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").Value = "288.81"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +2.30%  "

$ws.Range("D8").Value = "0.3284"
$ws.Range("E8").Value = "  -1.74%  "

$ws.Range("D9").Value = "44.46"
$ws.Range("E9").Value = "  -8.56%  "

$ws.Range("D10").Value = "1.151"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").Value = "0.07399"
$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("D12").Value = "0.9989"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").Value = "20.45"
$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").Value = "5.889"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "6.783"
$ws.Range("E15").Value = "  -2.58%  "

$ws.Range("D16").Value = "1.543.03"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "0.00001082"
$ws.Range("E17").Value = "  -3.23%  "

$ws.Range("D18").Value = "0.06652"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").Value = "86.37"
$ws.Range("E19").Value = "  -2.68%  "

$ws.Range("D20").Value = "6.432"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "0.9996"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").Value = "16.20"
$ws.Range("E22").Value = "  -2.26%  "

$ws.Range("D23").Value = "11.76"
$ws.Range("E23").Value = "  -3.56%  "

$ws.Range("D24").Value = "22.298.40"
$ws.Range("E24").Value = "  -1.13%  "

$ws.Range("D25").Value = "2.299"
$ws.Range("E25").Value = "  -4.26%  "

$ws.Range("D26").Value = "2.614"
$ws.Range("E26").Value = "  +0.33%  "

$ws.Range("D27").Value = "151.74"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("D28").Value = "19.41"
$ws.Range("E28").Value = "  -1.36%  "

$ws.Range("D29").Value = "4.937"
$ws.Range("E29").Value = "  -1.65%  "

$ws.Range("D30").Value = "123.17"
$ws.Range("E30").Value = "  -1.03%  "

$ws.Range("D31").Value = "1.720.82"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").Value = "1.095"
$ws.Range("E32").Value = "  +2.28%  "

$ws.Range("D33").Value = "5.957"
$ws.Range("E33").Value = "  -3.65%  "

$ws.Range("D34").Value = "1.910"
$ws.Range("E34").Value = "  -4.75%  "

$ws.Range("D35").Value = "9.470"
$ws.Range("E35").Value = "  -2.29%  "

$ws.Range("D36").Value = "0.08242"
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("D37").Value = "0.02375"
$ws.Range("E37").Value = "  -3.33%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "5.365"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.06335"
$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").Value = "0.2162"
$ws.Range("E40").Value = "  -4.73%  "

$ws.Range("D41").Value = "1.253"
$ws.Range("E41").Value = "  -3.79%  "

$ws.Range("D42").Value = "11.13"
$ws.Range("E42").Value = "  -2.38%  "

$ws.Range("D43").Value = "0.6099"
$ws.Range("E43").Value = "  -4.15%  "

$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "13.80"
$ws.Range("E45").Value = "  -1.47%  "

$ws.Range("D46").Value = "0.5955"
$ws.Range("E46").Value = "  -3.98%  "

$ws.Range("D47").Value = "3.749"
$ws.Range("E47").Value = "  -0.57%  "

$ws.Range("D48").Value = "1.994"
$ws.Range("E48").Value = "  -3.28%  "

$ws.Range("D49").Value = "123.35"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").Value = "1.181"
$ws.Range("E50").Value = "  -3.16%  "

$ws.Range("D51").Value = "0.07103"
$ws.Range("E51").Value = "  -2.29%  "

# Restore default (unstyled) formatting now that the text values are committed.
$priceRange.Style = "Normal"
